# Auto-generated edit script: bump the '想去人数' (F column, attendee-interest count)
# values across all four worksheets to match the upstream data refresh.
$wb = $excel.ActiveWorkbook
$appliedCount = 0
$mismatchCount = 0

# Sheet 1
$ws = $wb.Worksheets.Item(1)
$changes = @(
    @("F2", 41, 42),
    @("F4", 3505, 3519),
    @("F5", 3505, 3519),
    @("F6", 250, 251),
    @("F7", 5037, 5046),
    @("F9", 340, 341),
    @("F11", 675, 680),
    @("F13", 72, 74),
    @("F15", 687, 689),
    @("F16", 305, 306),
    @("F19", 154, 156),
    @("F22", 4870, 4877),
    @("F26", 5992, 5995),
    @("F28", 14, 15),
    @("F29", 3213, 3214),
    @("F30", 324, 325),
    @("F31", 696, 699),
    @("F33", 314, 316),
    @("F36", 971, 981),
    @("F40", 851, 852),
    @("F41", 945, 954)
)
foreach ($chg in $changes) {
    $ref = $chg[0]
    $oldVal = $chg[1]
    $newVal = $chg[2]
    $cell = $ws.Range($ref)
    $cur = $cell.Value2
    if ($cur -eq $oldVal) {
        $cell.Value2 = $newVal
        $appliedCount = $appliedCount + 1
    } else {
        $mismatchCount = $mismatchCount + 1
        Write-Output "Mismatch on sheet 1 ${ref}: expected $oldVal but found $cur"
    }
}

# Sheet 2
$ws = $wb.Worksheets.Item(2)
$changes = @(
    @("F3", 42, 44),
    @("F4", 22, 23)
)
foreach ($chg in $changes) {
    $ref = $chg[0]
    $oldVal = $chg[1]
    $newVal = $chg[2]
    $cell = $ws.Range($ref)
    $cur = $cell.Value2
    if ($cur -eq $oldVal) {
        $cell.Value2 = $newVal
        $appliedCount = $appliedCount + 1
    } else {
        $mismatchCount = $mismatchCount + 1
        Write-Output "Mismatch on sheet 2 ${ref}: expected $oldVal but found $cur"
    }
}

# Sheet 3
$ws = $wb.Worksheets.Item(3)
$changes = @(
    @("F2", 223, 224),
    @("F3", 1109, 1111)
)
foreach ($chg in $changes) {
    $ref = $chg[0]
    $oldVal = $chg[1]
    $newVal = $chg[2]
    $cell = $ws.Range($ref)
    $cur = $cell.Value2
    if ($cur -eq $oldVal) {
        $cell.Value2 = $newVal
        $appliedCount = $appliedCount + 1
    } else {
        $mismatchCount = $mismatchCount + 1
        Write-Output "Mismatch on sheet 3 ${ref}: expected $oldVal but found $cur"
    }
}

# Sheet 4
$ws = $wb.Worksheets.Item(4)
$changes = @(
    @("F2", 223, 224),
    @("F3", 41, 42),
    @("F4", 1109, 1111),
    @("F8", 3505, 3519),
    @("F9", 3505, 3519),
    @("F10", 250, 251),
    @("F11", 5037, 5046),
    @("F13", 340, 341),
    @("F15", 675, 680),
    @("F16", 72, 74),
    @("F18", 687, 689),
    @("F19", 305, 306),
    @("F21", 42, 44),
    @("F23", 154, 156),
    @("F26", 4870, 4877),
    @("F30", 5992, 5995),
    @("F32", 14, 15),
    @("F33", 3213, 3214),
    @("F34", 324, 325),
    @("F35", 696, 699),
    @("F37", 314, 316),
    @("F38", 22, 23),
    @("F41", 971, 981),
    @("F45", 851, 852),
    @("F46", 945, 954)
)
foreach ($chg in $changes) {
    $ref = $chg[0]
    $oldVal = $chg[1]
    $newVal = $chg[2]
    $cell = $ws.Range($ref)
    $cur = $cell.Value2
    if ($cur -eq $oldVal) {
        $cell.Value2 = $newVal
        $appliedCount = $appliedCount + 1
    } else {
        $mismatchCount = $mismatchCount + 1
        Write-Output "Mismatch on sheet 4 ${ref}: expected $oldVal but found $cur"
    }
}

Write-Output "Applied: $appliedCount, Mismatches: $mismatchCount"
